$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142330722250378
$ws.Range("D2").Value = 0.1135556752569897
$ws.Range("E2").Value = 0.130026972898321
$ws.Range("F2").Value = 2.004395868106442
$ws.Range("G2").Value = 1.354728967457035
$ws.Range("H2").Value = 1.255659114047816
$ws.Range("J2").Value = 0.1733118252852037
$ws.Range("K2").Value = 0.8694532474815162
$ws.Range("M2").Value = 0.3397596844345543
$ws.Range("N2").Value = 2.03996805465453
$ws.Range("B3").Value = 0.1329949848305745
$ws.Range("D3").Value = 0.1124055097550496
$ws.Range("E3").Value = 0.1296257360050959
$ws.Range("F3").Value = 1.999491145635005
$ws.Range("G3").Value = 1.348263359812805
$ws.Range("H3").Value = 1.258138806970464
$ws.Range("J3").Value = 0.1734199700604506
$ws.Range("K3").Value = 0.7935265334153314
$ws.Range("M3").Value = 0.3233113150382962
$ws.Range("N3").Value = 2.061534953335295
$ws.Range("B4").Value = 0.127333784714267
$ws.Range("D4").Value = 0.111735593348989
$ws.Range("E4").Value = 0.1294308185077817
$ws.Range("F4").Value = 1.997545158715226
$ws.Range("G4").Value = 1.345120191859323
$ws.Range("H4").Value = 1.260244265287895
$ws.Range("J4").Value = 0.1735624399010618
$ws.Range("K4").Value = 0.7471974914203372
$ws.Range("M4").Value = 0.3133712246041256
$ws.Range("N4").Value = 2.075451500095419
$ws.Range("B5").Value = 0.1250447898658393
$ws.Range("D5").Value = 0.1114717647386598
$ws.Range("E5").Value = 0.1293643455386011
$ws.Range("F5").Value = 1.997019964214857
$ws.Range("G5").Value = 1.344046917965386
$ws.Range("H5").Value = 1.261248826430275
$ws.Range("J5").Value = 0.1736396325105858
$ws.Range("K5").Value = 0.728391330008435
$ws.Range("M5").Value = 0.3093607097757527
$ws.Range("N5").Value = 2.081292136660785
$ws.Range("B6").Value = 0.1246657947446579
$ws.Range("D6").Value = 0.1114285111593176
$ws.Range("E6").Value = 0.129354090992134
$ws.Range("F6").Value = 1.996948927608258
$ws.Range("G6").Value = 1.343881231075983
$ws.Range("H6").Value = 1.261424485274475
$ws.Range("J6").Value = 0.1736536062017393
$ws.Range("K6").Value = 0.7252730176286377
$ws.Range("M6").Value = 0.3086971950103603
$ws.Range("N6").Value = 2.082272207786414
$ws.Range("B7").Value = 0.127302841499187
$ws.Range("D7").Value = 0.1117319980865759
$ws.Range("E7").Value = 0.1294298695366791
$ws.Range("F7").Value = 1.997536991573512
$ws.Range("G7").Value = 1.345104877176354
$ws.Range("H7").Value = 1.260257219712059
$ws.Range("J7").Value = 0.1735634034611451
$ws.Range("K7").Value = 0.7469435677510887
$ws.Range("M7").Value = 0.313316974685911
$ws.Range("N7").Value = 2.075529582608048
$ws.Range("B8").Value = 0.1390971097573299
$ws.Range("D8").Value = 0.113151591042822
$ws.Range("E8").Value = 0.129877961814632
$ws.Range("F8").Value = 2.002483485004177
$ws.Range("G8").Value = 1.352327849842908
$ws.Range("H8").Value = 1.256393119224683
$ws.Range("J8").Value = 0.1733333293213946
$ws.Range("K8").Value = 0.8432136310779015
$ws.Range("M8").Value = 0.334055284038925
$ws.Range("N8").Value = 2.047264325731838
$ws.Range("B9").Value = 0.1627839018027259
$ws.Range("D9").Value = 0.1162217221377944
$ws.Range("E9").Value = 0.1311640811612875
$ws.Range("F9").Value = 2.020647066212575
$ws.Range("G9").Value = 1.373068099246694
$ws.Range("H9").Value = 1.253442472480245
$ws.Range("J9").Value = 0.1734856029064673
$ws.Range("K9").Value = 1.034298854726671
$ws.Range("M9").Value = 0.3759848319681467
$ws.Range("N9").Value = 1.99718849142096
$ws.Range("B10").Value = 0.1805221992838284
$ws.Range("D10").Value = 0.1186500427163395
$ws.Range("E10").Value = 0.13235659142941
$ws.Range("F10").Value = 2.039168704083281
$ws.Range("G10").Value = 1.392341273475154
$ws.Range("H10").Value = 1.254099366195845
$ws.Range("J10").Value = 0.1739654671096531
$ws.Range("K10").Value = 1.176102284864044
$ws.Range("M10").Value = 0.4075606927868094
$ws.Range("N10").Value = 1.963661220770145
$ws.Range("B11").Value = 0.1886637947996093
$ws.Range("D11").Value = 0.1197918804822038
$ws.Range("E11").Value = 0.132952715596943
$ws.Range("F11").Value = 2.048722869191721
$ws.Range("G11").Value = 1.401991344945685
$ws.Range("H11").Value = 1.255012609960573
$ws.Range("J11").Value = 0.1742637219604077
$ws.Range("K11").Value = 1.240922553374673
$ws.Range("M11").Value = 0.4220930798760065
$ws.Range("N11").Value = 1.949117561993742
$ws.Range("B12").Value = 0.1917570912414845
$ws.Range("D12").Value = 0.1202295742363475
$ws.Range("E12").Value = 0.1331861496760425
$ws.Range("F12").Value = 2.052503320412882
$ws.Range("G12").Value = 1.405772906928121
$ws.Range("H12").Value = 1.255446844829095
$ws.Range("J12").Value = 0.1743881613290696
$ws.Range("K12").Value = 1.265513276157094
$ws.Range("M12").Value = 0.4276202867719263
$ws.Range("N12").Value = 1.943712139737897
$ws.Range("B13").Value = 0.1910904410440395
$ws.Range("D13").Value = 0.120135073803219
$ws.Range("E13").Value = 0.1331355335485433
$ws.Range("F13").Value = 2.051681902188349
$ws.Range("G13").Value = 1.404952812723167
$ws.Range("H13").Value = 1.255349391630034
$ws.Range("J13").Value = 0.1743608497897071
$ws.Range("K13").Value = 1.260215241622404
$ws.Range("M13").Value = 0.4264288333379085
$ws.Range("N13").Value = 1.944871758012578
$ws.Range("B14").Value = 0.1889180775966253
$ws.Range("D14").Value = 0.1198277837405186
$ws.Range("E14").Value = 0.1329717662540872
$ws.Range("F14").Value = 2.049030631319411
$ws.Range("G14").Value = 1.402299903121246
$ws.Range("H14").Value = 1.255046562539775
$ws.Range("J14").Value = 0.1742737292596885
$ws.Range("K14").Value = 1.242944753140591
$ws.Range("M14").Value = 0.4225473239593924
$ws.Range("N14").Value = 1.948670810768373
$ws.Range("B15").Value = 0.1875887726667713
$ws.Range("D15").Value = 0.1196402492672419
$ws.Range("E15").Value = 0.1328724555760523
$ws.Range("F15").Value = 2.047427819394542
$ws.Range("G15").Value = 1.400691507287092
$ws.Range("H15").Value = 1.254872586188526
$ws.Range("J15").Value = 0.1742218625853198
$ws.Range("K15").Value = 1.232371890702723
$ws.Range("M15").Value = 0.420172923802923
$ws.Range("N15").Value = 1.951011120442942
$ws.Range("B16").Value = 0.1799915697099692
$ws.Range("D16").Value = 0.1185761655246438
$ws.Range("E16").Value = 0.1323187111370387
$ws.Range("F16").Value = 2.038567044520079
$ws.Range("G16").Value = 1.391728410626342
$ws.Range("H16").Value = 1.25405205203964
$ws.Range("J16").Value = 0.1739475843442904
$ws.Range("K16").Value = 1.171872400576547
$ws.Range("M16").Value = 0.4066143462197189
$ws.Range("N16").Value = 1.964625941901875
$ws.Range("B17").Value = 0.1753493543899936
$ws.Range("D17").Value = 0.1179328778306754
$ws.Range("E17").Value = 0.1319927328942221
$ws.Range("F17").Value = 2.03342044923231
$ws.Range("G17").Value = 1.386456161642883
$ws.Range("H17").Value = 1.253706088015008
$ws.Range("J17").Value = 0.1737998024591292
$ws.Range("K17").Value = 1.134837891322462
$ws.Range("M17").Value = 0.3983396471325378
$ws.Range("N17").Value = 1.973159635291525
$ws.Range("B18").Value = 0.1726860945207562
$ws.Range("D18").Value = 0.1175663780732847
$ws.Range("E18").Value = 0.1318102893234325
$ws.Range("F18").Value = 2.030566478292187
$ws.Range("G18").Value = 1.383506756320259
$ws.Range("H18").Value = 1.253564927196976
$ws.Range("J18").Value = 0.1737223287406096
$ws.Range("K18").Value = 1.11356611889596
$ws.Range("M18").Value = 0.3935961129688224
$ws.Range("N18").Value = 1.978134628582461
$ws.Range("B19").Value = 0.1717855364396002
$ws.Range("D19").Value = 0.1174428903678901
$ws.Range("E19").Value = 0.1317493850835803
$ws.Range("F19").Value = 2.029618408948679
$ws.Range("G19").Value = 1.382522391494177
$ws.Range("H19").Value = 1.2535270629987
$ws.Range("J19").Value = 0.1736973902128653
$ws.Range("K19").Value = 1.106368938973418
$ws.Range("M19").Value = 0.3919927599173292
$ws.Range("N19").Value = 1.979830517272665
$ws.Range("B20").Value = 0.1758428214453573
$ws.Range("D20").Value = 0.1180009947352758
$ws.Range("E20").Value = 0.1320269112684542
$ws.Range("F20").Value = 2.033957318777311
$ws.Range("G20").Value = 1.387008803202633
$ws.Range("H20").Value = 1.253736931031199
$ws.Range("J20").Value = 0.173814755140782
$ws.Range("K20").Value = 1.138777227629987
$ws.Range("M20").Value = 0.3992188629157667
$ws.Range("N20").Value = 1.972244311307897
$ws.Range("B21").Value = 0.189555876418936
$ws.Range("D21").Value = 0.1199178986521616
$ws.Range("E21").Value = 0.1330196599969717
$ws.Range("F21").Value = 2.04980496206565
$ws.Range("G21").Value = 1.403075669409247
$ws.Range("H21").Value = 1.255133110920752
$ws.Range("J21").Value = 0.1742990066656347
$ws.Range("K21").Value = 1.248016303963084
$ws.Range("M21").Value = 0.4236867639112347
$ws.Range("N21").Value = 1.947552169076882
$ws.Range("B22").Value = 0.1985778292460338
$ws.Range("D22").Value = 0.1212016058559158
$ws.Range("E22").Value = 0.1337133209765611
$ws.Range("F22").Value = 2.061109564481725
$ws.Range("G22").Value = 1.414318409259835
$ws.Range("H22").Value = 1.256560939879023
$ws.Range("J22").Value = 0.1746825041276381
$ws.Range("K22").Value = 1.319670710598302
$ws.Range("M22").Value = 0.4398184200021831
$ws.Range("N22").Value = 1.932008705778387
$ws.Range("B23").Value = 0.1937572334660445
$ws.Range("D23").Value = 0.1205136537285298
$ws.Range("E23").Value = 0.1333390043278442
$ws.Range("F23").Value = 2.054989342275761
$ws.Range("G23").Value = 1.408249920537145
$ws.Range("H23").Value = 1.255751707393529
$ws.Range("J23").Value = 0.1744716933098829
$ws.Range("K23").Value = 1.281403694047924
$ws.Range("M23").Value = 0.4311958327555772
$ws.Range("N23").Value = 1.9402501252024
$ws.Range("B24").Value = 0.1756197074993793
$ws.Range("D24").Value = 0.117970188692432
$ws.Range("E24").Value = 0.1320114437569266
$ws.Range("F24").Value = 2.033714273381847
$ws.Range("G24").Value = 1.386758699544089
$ws.Range("H24").Value = 1.253722807045591
$ws.Range("J24").Value = 0.1738079717090457
$ws.Range("K24").Value = 1.136996191799767
$ws.Range("M24").Value = 0.3988213267771243
$ws.Range("N24").Value = 1.972657914786215
$ws.Range("B25").Value = 0.1563167180824081
$ws.Range("D25").Value = 0.1153607105735617
$ws.Range("E25").Value = 0.1307726158642843
$ws.Range("F25").Value = 2.014825486329386
$ws.Range("G25").Value = 1.366750585672818
$ws.Range("H25").Value = 1.253744951177183
$ws.Range("J25").Value = 0.1733797950247009
$ws.Range("K25").Value = 0.9823578942436484
$ws.Range("M25").Value = 0.3645065962029861
$ws.Range("N25").Value = 2.010162292433145
